$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Align F7:G7 formatting with F8's style (removes the now-unused duplicate style)
$ws.Range("F8").Copy() | Out-Null
$ws.Range("F7:G7").PasteSpecial(-4122) | Out-Null

# Mark "Crear el proyecto en Django" task as completed
$ws.Range("D7").Value = 1

# Mark "Crear los templates..." task as completed and fix typo in its note
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = "Está bastante avanzado"

# Fill in the new task row (row 9)
$ws.Range("B9").Value = "Crear la segunda parte con la consulta de los datos de las pelis"
$ws.Range("C9").Value = "Alta"
$ws.Range("E9").Value = "Probablemente se pueda hacer directamente con una API call similar, porque trae mucha más información, hasta imágenes"

# Update selection to match the author's last position
$ws.Range("F11").Select() | Out-Null
